$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) "Highly inclined towards Cloud and Vue js framework. Researched
#    integration of Cloud platforms like Firebase & AWS to Vue Js and
#    Azure to .Net core." paragraph gets reworded to also call out
#    React js, in three separate Find/Replace passes that line up with
#    the untouched "js" run (wrapped in proofErr/spellStart-spellEnd)
#    sitting between them.
# ----------------------------------------------------------------------

$rng = $d.Content
$found1 = $rng.Find.Execute(
    "Highly inclined towards Cloud and Vue ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Highly inclined towards Cloud applications, Vue ", 2)
Write-Host "Step1a replaced: $found1"

$rng = $d.Content
$found2 = $rng.Find.Execute(
    " framework. Researched integration of Cloud platforms like Firebase & AWS to Vue ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " & React js frameworks. Researched integration of Cloud platforms like Firebase & AWS to Vue ", 2)
Write-Host "Step1b replaced: $found2"

$rng = $d.Content
$found3 = $rng.Find.Execute(
    "Js and Azure", $true, $false, $false, $false, $false, $true, 1, $false,
    "js, React js and Azure", 2)
Write-Host "Step1c replaced: $found3"

# ----------------------------------------------------------------------
# 2) Skills table cell: "Vue.js, Angular " -> "Vue.js, React.js, Angular "
#    (scope the search to the table so the similar-looking phrase higher
#    up, in the "Experience in front-end frameworks like ..." sentence,
#    is left untouched).
# ----------------------------------------------------------------------

$tblRng = $d.Tables(1).Range
$found4 = $tblRng.Find.Execute(
    "Vue.js, Angular ", $true, $false, $false, $false, $false, $true, 1,
    $false, "Vue.js, React.js, Angular ", 2)
Write-Host "Step2 replaced: $found4"

# ----------------------------------------------------------------------
# 3) The stock-tracking bullet becomes a description of the React chat
#    project.
# ----------------------------------------------------------------------

$rng = $d.Content
$found5 = $rng.Find.Execute(
    "A Vue project with Firebase authentication view tracked stocks with graph",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A React PWA single room chat project with Firebase google authentication and Firestore data for real time listener",
    2)
Write-Host "Step3 replaced: $found5"

# ----------------------------------------------------------------------
# 4) The plain-text stock-app URL underneath becomes a live hyperlink
#    pointing at the new React chat app.
# ----------------------------------------------------------------------

$rng = $d.Content
$found6 = $rng.Find.Execute(
    "https://stockapi-90a27.firebaseapp.com", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0)
Write-Host "Step4 found old link: $found6"
if ($found6) {
    $rng.Text = ""
    $hyperlink = $d.Hyperlinks.Add(
        $rng,
        "https://react-superchat-8b806.web.app/",
        [System.Reflection.Missing]::Value,
        [System.Reflection.Missing]::Value,
        "https://react-superchat-8b806.web.app/")
    Write-Host "Step4 new hyperlink text: $($hyperlink.Range.Text)"
}
